$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(44, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/01"
$dateCell.ClearFormats()

$ws.Cells.Item(44, 2).Value = "水"
$ws.Cells.Item(44, 3).Value = 9
$ws.Cells.Item(44, 4).Value = 3
